# Apply coinranking snapshot refresh: update Price (D) and Volume(1h) (E) columns
# for the updated rows per the Sun Feb 11 05:44:05 UTC 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "48.268.18"
$ws.Range("E2").Value = "  +1.82%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.522.93"
$ws.Range("E3").Value = "  +0.77%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'323.45"

# Row 6 - Solana
$ws.Range("D6").Value = "'109.02"
$ws.Range("E6").Value = "  -0.96%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.527"
$ws.Range("E7").Value = "  +0.46%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.01%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.561"
$ws.Range("E9").Value = "  +4.29%  "

# Row 10 - Avalanche
$ws.Range("E10").Value = "  +2.26%  "

# Row 11 - Chainlink
$ws.Range("D11").Value = "'20.25"
$ws.Range("E11").Value = "  +9.39%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  +0.44%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.86%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'7.27"
$ws.Range("E14").Value = "  +0.73%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.916.84"
$ws.Range("E15").Value = "  +0.79%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.514.31"
$ws.Range("E16").Value = "  +0.52%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "'0.862"
$ws.Range("E17").Value = "  +0.15%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "48.150.74"
$ws.Range("E18").Value = "  +1.72%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Range("D19").Value = "'13.23"
$ws.Range("E19").Value = "  +2.70%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  -0.37%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0₃0945"
$ws.Range("E21").Value = "  +0.39%  "

# Row 22 - ImmutableX
$ws.Range("E22").Value = "  +0.47%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "'72.44"
$ws.Range("E23").Value = "  +2.65%  "

# Row 24 - BitcoinCash
$ws.Range("D24").Value = "'267.83"
$ws.Range("E24").Value = "  +7.75%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "'2.57"
$ws.Range("E25").Value = "  -1.32%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "'26.21"
$ws.Range("E26").Value = "  +0.35%  "

# Row 29 - Cosmos
$ws.Range("E29").Value = "  +0.93%  "

# Row 30 - Kaspa
$ws.Range("E30").Value = "  +4.87%  "

# Row 31 - InjectiveProtocol
$ws.Range("D31").Value = "'34.97"
$ws.Range("E31").Value = "  -1.10%  "

# Row 32 - OKB
$ws.Range("D32").Value = "'49.80"
$ws.Range("E32").Value = "  -0.29%  "

# Row 33 - Celestia
$ws.Range("D33").Value = "'20.01"
$ws.Range("E33").Value = "  +0.02%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "'5.39"
$ws.Range("E34").Value = "  -0.97%  "

# Row 35 - FirstDigitalUSD
$ws.Range("E35").Value = "  -0.02%  "

# Row 36 - Hedera
$ws.Range("D36").Value = "'0.0793"
$ws.Range("E36").Value = "  -0.54%  "

# Row 37 - ARBITRUM
$ws.Range("E37").Value = "  -0.34%  "

# Row 38 - RenderToken
$ws.Range("D38").Value = "'4.73"
$ws.Range("E38").Value = "  +0.48%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  -0.32%  "

# Row 40 - Stellar
$ws.Range("E40").Value = "  +0.15%  "

# Row 41 - EnergySwap
$ws.Range("D41").Value = "'22.14"
$ws.Range("E41").Value = "  +4.22%  "

# Row 42 - WEMIXToken
$ws.Range("E42").Value = "  -1.46%  "

# Row 43 - Monero
$ws.Range("D43").Value = "'119.04"
$ws.Range("E43").Value = "  -2.21%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  +0.00%  "

# Row 45 - Maker
$ws.Range("D45").Value = "2.002.80"
$ws.Range("E45").Value = "  +0.08%  "

# Row 46 - NEARProtocol
$ws.Range("E46").Value = "  +0.88%  "

# Row 47 - Stacks
$ws.Range("D47").Value = "'1.89"
$ws.Range("E47").Value = "  +6.22%  "

# Row 48 - ApeXProtocol
$ws.Range("E48").Value = "  -1.77%  "

# Row 49 - FraxShare
$ws.Range("D49").Value = "'9.09"
$ws.Range("E49").Value = "  +0.27%  "

# Row 50 - THORChain
$ws.Range("D50").Value = "'5.27"
$ws.Range("E50").Value = "  +0.55%  "

# Row 51 - BitcoinSV
$ws.Range("D51").Value = "'80.58"
$ws.Range("E51").Value = "  +3.02%  "
